$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the whole existing data block (A2:D15) - the "testing script" reshuffled
# / re-sampled the set of teams it reports on, so start from a clean slate before
# writing the new rows back out.
$ws.Range("A2:D15").ClearContents()

# New team order (rows 2-13). Only Minnesota Twins / Texas Rangers / Pittsburgh
# Pirates carry weather+odds data this run - the rest are placeholders with just
# the team name, same as the other untouched rows in the sheet.
$ws.Range("A2").Value = "Atlanta Braves"
$ws.Range("A3").Value = "Toronto Blue Jays"
$ws.Range("A4").Value = "St. Louis Cardinals"
$ws.Range("A5").Value = "Baltimore Orioles"
$ws.Range("A6").Value = "Cleveland Indians"

$ws.Range("A7").Value = "Minnesota Twins"
$ws.Range("B7").Value = 64
$ws.Range("C7").Value = 102
$ws.Range("D7").Value = 8

$ws.Range("A8").Value = "Texas Rangers"
$ws.Range("B8").Value = 79
$ws.Range("C8").Value = -172
$ws.Range("D8").Value = 9.5

$ws.Range("A9").Value = "Pittsburgh Pirates"
$ws.Range("B9").Value = 64
$ws.Range("C9").Value = -174
$ws.Range("D9").Value = 8.5

$ws.Range("A10").Value = "Los Angeles Dodgers"
$ws.Range("A11").Value = "Houston Astros"
$ws.Range("A12").Value = "Arizona Diamondbacks"
$ws.Range("A13").Value = "Chicago White Sox"

# Selection moved to B8 by the time the script finished running.
$ws.Range("B8").Select()

# The workbook no longer carries any Hyperlink / Followed Hyperlink cell
# styles (they were unused leftovers).
try { $wb.Styles.Item("Hyperlink").Delete() } catch {}
try { $wb.Styles.Item("Followed Hyperlink").Delete() } catch {}
